# OpenTBS demo doc: arrange examples and doc
#
# 1) Split the "[onshow..now;frm='yyyy-mm-dd hh:nn:ss']" run into five runs
#    (one per logical piece) and normalise the curly quotes around the
#    date-format string to straight apostrophes.
# 2) Move the "_GoBack" bookmark from the end of the chart paragraph to
#    right after the "[onshow..cst.PHP_VERSION]" run.

$d = $word.ActiveDocument

# --- 1) split/normalise the date-format field -----------------------------
$rng = $d.Content
$found = $rng.Find.Execute("[onshow..now;frm=’yyyy-mm-dd hh:nn:ss’]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $rng.Start

    # Boundaries (relative to $start) of the five pieces:
    #   [0, 17)  -> "[onshow..now;frm="
    #   [17,18)  -> opening quote
    #   [18,37)  -> "yyyy-mm-dd hh:nn:ss"
    #   [37,38)  -> closing quote
    #   [38,39)  -> "]"
    $b1 = $start + 17
    $b2 = $start + 18
    $b3 = $start + 37
    $b4 = $start + 38
    $b5 = $start + 39

    # Replace the curly quotes with straight apostrophes first (while the
    # whole field is still one run) ...
    $q1 = $d.Range($b1, $b2)
    $q1.Text = "'"
    $q2 = $d.Range($b3, $b4)
    $q2.Text = "'"

    # ... then force the run to split at each internal boundary by toggling
    # a character property on/off (the net formatting doesn't change, but
    # it creates a fresh run boundary, matching how Word itself breaks runs
    # when content is edited piecewise).
    $r1 = $d.Range($b1, $b2)
    $r1.Bold = 1
    $r1.Bold = 0

    $r2 = $d.Range($b2, $b3)
    $r2.Bold = 1
    $r2.Bold = 0

    $r3 = $d.Range($b3, $b4)
    $r3.Bold = 1
    $r3.Bold = 0

    $r4 = $d.Range($b4, $b5)
    $r4.Bold = 1
    $r4.Bold = 0
}

# --- 2) relocate the "_GoBack" bookmark ------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("[onshow..cst.PHP_VERSION]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    $insertAt = $rng2.End
    $newBmRange = $d.Range($insertAt, $insertAt)
    $d.Bookmarks.Add("_GoBack", $newBmRange)
}
